$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note about menu temp art, added in row 21 (existing blank row in the
# "NIFTY SHIT TO ADD" block) - set first so it lands as shared string #48.
$ws.Range("C21").Value = "Menus need to at the very least use better temp art"

# New "Level Designers" person tag next to the "maps need balancing" line
# (row 15) - set second so it lands as shared string #49.
$ws.Range("A15").Value = "Level Designers"

# Insert a new blank row at 23 (pushes the "NIFTY SHIT TO ADD" section and
# everything below it down by one row).
$ws.Rows(23).Insert()

# Update the view so the window is scrolled/selected around row 16, matching
# the author's new cursor position.
$ws.Range("A16").Select()
